# Update numeric result values (columns B:H, rows 2-11) on the Results sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Results")

$data = @{
    2  = @(53.42058726655018, 48.21134383762951, 15.8649639555416, 215.113685541457, 185.6645927840595, -2.128363103746001, 0.9685436417123128)
    3  = @(93.92634631660378, 32.57606738439313, 21.62976589917799, 217.9416862986527, 177.7674827192845, 0.0266539393121521, 0.8434430846515169)
    4  = @(42.07872482300356, 55.778753620023, 13.75665567428846, 212.4042245709743, 165.9930991295923, 0.6606014699702949, 0.628580577673582)
    5  = @(25.1320689820177, 53.95442736672652, 23.37174810531702, 220.7258896500101, 168.8606116161002, -2.455548055447991, 0.9822875583611832)
    6  = @(87.83518319932872, 60.85261565329055, 42.91725817435447, 209.0091101413919, 169.7765287996725, -0.2690271177533865, 0.6124472030349173)
    7  = @(88.98508311135974, 67.41360589240355, 13.99778848415658, 229.5392251809899, 193.5399550936415, 0.4087644952283602, 0.8054540780090966)
    8  = @(23.21790560505868, 67.19586573437418, 14.84825188608004, 212.166344470336, 165.3424558863049, -2.433976775792443, 0.96236875154338)
    9  = @(43.63248597736443, 56.45469626035633, 10.76868152101614, 219.4419299805639, 183.5252236095589, -2.74875737596758, 0.6208635236352059)
    10 = @(71.8854858177927, 19.10522070429979, 31.50592525958717, 203.7675879559816, 152.2825251495386, -1.461363153723463, 0.7520746459758204)
    11 = @(63.28342874040409, 79.70305340599327, 45.67323006493222, 208.6646039776573, 178.2486916864511, -2.709647608058247, 0.7908960913419094)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($col = 2; $col -le 8; $col++) {
        $ws.Cells.Item($row, $col).Value = $values[$col - 2]
    }
}

$wb.Save()
